# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 1068, pushing the existing
# rows 1068-1104 down to 1069-1105 (dimension grows from T1104 to T1105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 1068. Excel shifts every row at/after
# 1068 down by one, and the new blank row 1068 inherits formatting (incl. the
# date number format in column D) from the row above, same as typical manual
# "Insert Row" behaviour.
$ws.Rows.Item(1068).Insert()

# Populate the freshly inserted row 1068 with the new record's values.
$ws.Cells.Item(1068, 1).Value2 = 9
$ws.Cells.Item(1068, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1068, 3).Value2 = "Metropolitana"
$ws.Cells.Item(1068, 4).Value2 = 45075
$ws.Cells.Item(1068, 5).Value2 = 13
$ws.Cells.Item(1068, 6).Value2 = "Fruta"
$ws.Cells.Item(1068, 7).Value2 = 100104
$ws.Cells.Item(1068, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(1068, 9).Value2 = 100104005
$ws.Cells.Item(1068, 10).Value2 = "Pera"
$ws.Cells.Item(1068, 11).Value2 = "Packham's Triumph"
$ws.Cells.Item(1068, 12).Value2 = "Primera"
$ws.Cells.Item(1068, 13).Value2 = 380
$ws.Cells.Item(1068, 14).Value2 = 11500
$ws.Cells.Item(1068, 15).Value2 = 12000
$ws.Cells.Item(1068, 16).Value2 = 11763
$ws.Cells.Item(1068, 17).Value2 = "$/caja 18 kilos granel"
$ws.Cells.Item(1068, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(1068, 19).Value2 = 654
$ws.Cells.Item(1068, 20).Value2 = 18
